$d = $word.ActiveDocument

# 1. Update portfolio URL
$d.Content.Find.Execute("www.jedulan.cf", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "www.portfolio-jedulan.netlify.app", 2) | Out-Null

# 2. Remove the _GoBack bookmark from its current location (end of the
#    "HTML5/ HTML" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the "TECHNICAL SKILLS" list items by their current text so the
# script is resilient to any paragraph-numbering quirks.
function Get-ParaByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

# 3. "Mysql, VB" -> "Wordpress" (keep the spell-check wrap around the
#    single remaining word, drop the ", VB" run entirely).
$pMysql = Get-ParaByText $d "Mysql, VB"
$findRange = $pMysql.Range.Duplicate
$findRange.Find.Execute("Mysql", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$findRange.Text = "Wordpress"

$vbRange = $pMysql.Range.Duplicate
$vbRange.Find.Execute(", VB", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$vbRange.Delete()

# 4. "Wordpress" -> "JSON", and drop its spell-check wrap entirely (the
#    new word replaces both the proofErr markers and the run).
$pWordpress = Get-ParaByText $d "Wordpress"
$pPrXml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr>'
$runXml = '<w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>JSON</w:t></w:r>'
$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrXml + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pWordpress.Range.InsertXML($packageXml)

# 5. "JSON" -> "56 wpm"
$pJson = Get-ParaByText $d "JSON"
$pJson.Range.Text = "56 wpm"

# 6. "56 wpm" -> "MySQL"
$pWpm = Get-ParaByText $d "56 wpm"
$pWpm.Range.Text = "MySQL"

# 7/8. Remove the (now duplicate) "MySQL" and "MySQL Server" paragraphs
#      entirely.
$pOldMysql = Get-ParaByText $d "MySQL"
$pOldMysql.Range.Delete()
$pOldMysqlServer = Get-ParaByText $d "MySQL Server"
$pOldMysqlServer.Range.Delete()

# 9. Re-insert the _GoBack bookmark, now collapsed at the start of the
#    "ReactJS" paragraph.
$pReact = Get-ParaByText $d "ReactJS"
$bmRange = $d.Range($pReact.Range.Start, $pReact.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
